$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.999898374080658
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.218034029006958
